# feat: add 2022-Q3 data
#
# Original workbook has two sheets:
#   1) "总计"      - a summary sheet, one data row for "2021-Q3"
#   2) "2021-Q3"   - per-fund holdings for the 2021-Q3 quarter
#
# The edit adds a new "2022-Q3" quarter:
#   - "总计" gets a new row on top with the 2022-Q3 totals, the old
#     2021-Q3 row is pushed down.
#   - A new sheet "2022-Q3" is created (reusing the slot/position of the
#     old "2021-Q3" sheet) holding the 2022-Q3 per-fund holdings.
#   - The original "2021-Q3" per-fund holdings are preserved, unchanged,
#     in a sheet of their own placed right after "2022-Q3".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "总计"
$ws2 = $wb.Worksheets.Item(2)   # "2021-Q3" (existing per-fund sheet)

# --------------------------------------------------------------------
# 1) Duplicate the existing "2021-Q3" sheet so its data is preserved
#    verbatim in its own sheet, placed right after the original.
# --------------------------------------------------------------------
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "2021-Q3-holder"

# Rename the original sheet - it will now hold the new 2022-Q3 data.
$ws2.Name = "2022-Q3"
$ws3.Name = "2021-Q3"

# --------------------------------------------------------------------
# 2) "总计" sheet: insert the 2022-Q3 totals as the new row 2, pushing
#    the existing 2021-Q3 row down to row 3 (keeping its formatting).
# --------------------------------------------------------------------
$ws1.Range("A2:D2").Copy($ws1.Range("A3:D3"))
$ws1.Range("A3").Value = 1

$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = 0.11

# --------------------------------------------------------------------
# 3) Rebuild the "2022-Q3" sheet (formerly the "2021-Q3" per-fund
#    sheet) with the new quarter's fund holdings.
# --------------------------------------------------------------------
$ws2.Cells.Clear()

# Header row - copy the header style used on "总计" (bold, bordered,
# centered) onto the whole header row, then fill in the labels.
$ws1.Range("B1").Copy($ws2.Range("B1:H1"))
$ws2.Range("B1").Value = "基金代码"
$ws2.Range("C1").Value = "基金名称"
$ws2.Range("D1").Value = "基金规模"
$ws2.Range("E1").Value = "股票总仓位"
$ws2.Range("F1").Value = "仓位占比"
$ws2.Range("G1").Value = "持有市值(亿元)"
$ws2.Range("H1").Value = "仓位排名"

# Column A (row index) - copy the same styled cell used on "总计".
$ws1.Range("A2").Copy($ws2.Range("A2:A3"))
$ws2.Range("A2").Value = 0
$ws2.Range("A3").Value = 1

# Row 2 - 005014 / 泰康景泰回报混合A
$ws2.Range("B2").Value = "'005014"
$ws2.Range("C2").Value = "泰康景泰回报混合A"
$ws2.Range("D2").Value = "'8.99"
$ws2.Range("E2").Value = "'34.29"
$ws2.Range("F2").Value = "'1.20"
$ws2.Range("G2").Value = "'0.1079"
$ws2.Range("H2").Value = 10

# Row 3 - 005015 / 泰康景泰回报混合C
$ws2.Range("B3").Value = "'005015"
$ws2.Range("C3").Value = "泰康景泰回报混合C"
$ws2.Range("D3").Value = "'0.39"
$ws2.Range("E3").Value = "'34.29"
$ws2.Range("F3").Value = "'1.20"
$ws2.Range("G3").Value = "'0.0047"
$ws2.Range("H3").Value = 10

# --------------------------------------------------------------------
# 4) Keep "总计" as the active sheet, as it was originally.
# --------------------------------------------------------------------
$ws1.Activate()
